# Updated symbol list with refreshed price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'317.23"
$ws.Range("E2").Value = "'2.11%"
$ws.Range("D3").Value = "'40.93"
$ws.Range("E3").Value = "'-1.33%"
$ws.Range("E4").Value = "'0.33%"
$ws.Range("D5").Value = "'0.07635"
$ws.Range("E5").Value = "'-0.83%"
$ws.Range("D6").Value = "'1.689"
$ws.Range("E6").Value = "'4.10%"
$ws.Range("D7").Value = "'0.9355"
$ws.Range("E7").Value = "'1.30%"
$ws.Range("E8").Value = "'-1.74%"
$ws.Range("D9").Value = "'0.1244"
$ws.Range("E9").Value = "'1.93%"
$ws.Range("D10").Value = "'0.1820"
$ws.Range("E10").Value = "'-1.11%"
$ws.Range("D11").Value = "'0.09030"
$ws.Range("E11").Value = "'-1.85%"
$ws.Range("D12").Value = "'0.04157"
$ws.Range("E12").Value = "'-3.83%"
$ws.Range("E13").Value = "'0.72%"
$ws.Range("D14").Value = "'0.001269"
$ws.Range("E14").Value = "'1.92%"
$ws.Range("D15").Value = "'0.005882"
$ws.Range("E15").Value = "'1.05%"
$ws.Range("D17").Value = "'3.356"
$ws.Range("D18").Value = "'4.333"
$ws.Range("E18").Value = "'0.54%"
$ws.Range("D19").Value = "'0.3359"
$ws.Range("E19").Value = "'1.67%"
$ws.Range("D20").Value = "'8.387"
$ws.Range("E20").Value = "'21.12%"
$ws.Range("D21").Value = "'0.1348"
$ws.Range("E21").Value = "'-2.99%"
$ws.Range("D22").Value = "'0.2739"
$ws.Range("E22").Value = "'2.28%"
$ws.Range("D23").Value = "'0.04041"
$ws.Range("E23").Value = "'-0.32%"
$ws.Range("D24").Value = "'0.001266"
$ws.Range("D25").Value = "'0.004085"
$ws.Range("E25").Value = "'-0.06%"
$ws.Range("D26").Value = "'0.0001275"
$ws.Range("E26").Value = "'0.54%"
$ws.Range("D38").Value = "'0.02503"
$ws.Range("E38").Value = "'1.46%"
$ws.Range("D39").Value = "'0.05203"
$ws.Range("E39").Value = "'-1.13%"
$ws.Range("D40").Value = "'0.007798"
$ws.Range("E40").Value = "'-0.34%"
$ws.Range("D41").Value = "'0.1300"
$ws.Range("E41").Value = "'-1.13%"
$ws.Range("D42").Value = "'0.007375"
$ws.Range("E42").Value = "'8.60%"
$ws.Range("D43").Value = "'0.002169"
$ws.Range("E43").Value = "'16.78%"
$ws.Range("D44").Value = "'0.008230"
$ws.Range("E44").Value = "'0.51%"
$ws.Range("D45").Value = "'0.3167"
$ws.Range("E45").Value = "'2.09%"
$ws.Range("D46").Value = "'0.00006657"
$ws.Range("E46").Value = "'-1.05%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.57%"
$ws.Range("D48").Value = "'0.2766"
$ws.Range("E48").Value = "'62.81%"
$ws.Range("D49").Value = "'0.004218"
$ws.Range("E49").Value = "'3.02%"
$ws.Range("D50").Value = "'0.00002109"
$ws.Range("E50").Value = "'0.57%"
$ws.Range("D51").Value = "'0.0002009"
$ws.Range("E51").Value = "'0.57%"
